$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" date field text that lives on
#    the slide master and every slide layout (Insert > Header & Footer >
#    Apply to All re-stamps this cached text with the current date).
#    12/2/2018 -> 12/4/2018
# ---------------------------------------------------------------------
function Find-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDatePh = $true }
        } catch {
        }
        if ($isDatePh) { return $sh }
    }
    return $null
}

$newDate = "12/4/2018"

$masterDateShape = Find-DatePlaceholder $p.SlideMaster
if ($masterDateShape -ne $null) {
    if ($masterDateShape.TextFrame.TextRange.Text -ne $newDate) {
        $masterDateShape.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    $layoutDateShape = Find-DatePlaceholder $layout
    if ($layoutDateShape -ne $null) {
        if ($layoutDateShape.TextFrame.TextRange.Text -ne $newDate) {
            $layoutDateShape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 10 ("log of literate rate vs. Log of survival rate"):
#    update the regression-summary text box.
# ---------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$slide10Body = $slide10.Shapes.Item(3)
$slide10Range = $slide10Body.TextFrame.TextRange

$slide10Range.Paragraphs(1, 1).Text = "Approximately 50% of the variation in the log of survival rate is explained by the log of the literate rate."
$slide10Range.Paragraphs(2, 1).Text = "It's a log-log transformation, so it resulted in an interpretation that a doubling of the literacy rate of youth results in an 85% increase in the rate of survival to the age of 65."

# ---------------------------------------------------------------------
# 3) Slide 11 ("Multicollinearity plots matrix"):
#    update the caption describing the matrix plot.
# ---------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)
$slide11Body = $slide11.Shapes.Item(3)
$slide11Body.TextFrame.TextRange.Text = "There is minimum relationship between mother’s position maternity, female life expectancy, female participation, and female manager."
